# Weekly CompStat update: new crime data collected.
# Updates the report header (volume/week-of text), and refreshes the
# weekly/28-day/YTD/2-year crime-complaint figures for the 13th Precinct.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text (rich-text shared strings): bump the report volume number
# and roll the "week covering" date range forward by one week.
# ---------------------------------------------------------------------
$ws.Range("A8").Value  = "Volume 32   Number  26"
$ws.Range("C9").Value  = "Report Covering the Week  6/23/2025  Through  6/29/2025"

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 6
$ws.Range("K15").Value = 83.333333333333

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 66.666666666666
$ws.Range("F16").Value = 18
$ws.Range("H16").Value = 63.636363636363
$ws.Range("I16").Value = 96
$ws.Range("J16").Value = 72
$ws.Range("K16").Value = 33.333333333333
$ws.Range("L16").Value = 20
$ws.Range("M16").Value = 20
$ws.Range("N16").Value = -81.573896353167

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -12.5
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = 38.095238095238
$ws.Range("I17").Value = 140
$ws.Range("J17").Value = 107
$ws.Range("K17").Value = 30.841121495327
$ws.Range("L17").Value = 33.333333333333
$ws.Range("M17").Value = 129.508196721311
$ws.Range("N17").Value = 6.870229007633

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 21
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 61.538461538461
$ws.Range("I18").Value = 172
$ws.Range("J18").Value = 92
$ws.Range("K18").Value = 86.956521739130
$ws.Range("L18").Value = 48.275862068965
$ws.Range("M18").Value = 24.637681159420
$ws.Range("N18").Value = -79.075425790754

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 22
$ws.Range("D19").Value = 25
$ws.Range("E19").Value = -12
$ws.Range("F19").Value = 76
$ws.Range("G19").Value = 79
$ws.Range("H19").Value = -3.797468354430
$ws.Range("I19").Value = 531
$ws.Range("J19").Value = 480
$ws.Range("K19").Value = 10.625
$ws.Range("L19").Value = 5.566600397614
$ws.Range("M19").Value = -24.466571834992
$ws.Range("N19").Value = -59.090909090909

# ---------------------------------------------------------------------
# Row 20 - G.L.A.  (Week-to-Date complaint count fell to zero, so C20
# switches from a numeric 1 to the shared "0" text used elsewhere in
# the sheet for zero counts.)
# ---------------------------------------------------------------------
$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4122)          # xlPasteFormats - adopt the "text zero" style
$ws.Range("C20").Formula = "=""0"""
$ws.Range("C20").Copy()
$ws.Range("C20").PasteSpecial(-4163)          # xlPasteValues - freeze the literal text "0"
$excel.CutCopyMode = $false

$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 5
$ws.Range("H20").Value = 66.666666666666
$ws.Range("J20").Value = 17
$ws.Range("K20").Value = 35.294117647058
$ws.Range("L20").Value = -32.352941176470
$ws.Range("N20").Value = -96.192052980132

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 36
$ws.Range("D21").Value = 39
$ws.Range("E21").Value = -7.692307692307
$ws.Range("F21").Value = 149
$ws.Range("G21").Value = 129
$ws.Range("H21").Value = 15.503875968992
$ws.Range("I21").Value = 973
$ws.Range("J21").Value = 775
$ws.Range("K21").Value = 25.548387096774
$ws.Range("L21").Value = 15.147928994082
$ws.Range("M21").Value = -2.991026919242
$ws.Range("N21").Value = -71.280991735537

# ---------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------
$ws.Range("C22").Value = 2
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 4
$ws.Range("I22").Value = 36
$ws.Range("J22").Value = 30
$ws.Range("K22").Value = 20
$ws.Range("L22").Value = -14.285714285714
$ws.Range("M22").Value = 0

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 41
$ws.Range("D24").Value = 52
$ws.Range("E24").Value = -21.153846153846
$ws.Range("F24").Value = 227
$ws.Range("G24").Value = 205
$ws.Range("H24").Value = 10.731707317073
$ws.Range("I24").Value = 1351
$ws.Range("J24").Value = 1455
$ws.Range("K24").Value = -7.147766323024
$ws.Range("L24").Value = 29.406130268199
$ws.Range("M24").Value = 54.223744292237

# ---------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 35
$ws.Range("D25").Value = 46
$ws.Range("E25").Value = -23.913043478260
$ws.Range("F25").Value = 171
$ws.Range("G25").Value = 174
$ws.Range("H25").Value = -1.724137931034
$ws.Range("I25").Value = 1086
$ws.Range("J25").Value = 1246
$ws.Range("K25").Value = -12.841091492776
$ws.Range("L25").Value = 41.960784313725

# ---------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C26").Value = 17
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = 54.545454545454
$ws.Range("F26").Value = 48
$ws.Range("G26").Value = 60
$ws.Range("H26").Value = -20
$ws.Range("I26").Value = 285
$ws.Range("J26").Value = 297
$ws.Range("K26").Value = -4.040404040404
$ws.Range("L26").Value = 20.762711864406
$ws.Range("M26").Value = 43.21608040201

# ---------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------
$ws.Range("D27").Value = 2
$ws.Range("G27").Value = 5
$ws.Range("J27").Value = 11
$ws.Range("K27").Value = 36.363636363636

# ---------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -50
$ws.Range("F28").Value = 8
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 33.333333333333
$ws.Range("I28").Value = 60
$ws.Range("J28").Value = 56
$ws.Range("K28").Value = 7.142857142857
$ws.Range("L28").Value = -3.225806451612

# ---------------------------------------------------------------------
# Row 31 - Hate Crimes (week-to-date went from no data to one reported,
# then back out, so D31/G31/E31/H31 flip from the shared "0"/"***.*"
# text placeholders to real numeric counts/percentages.)
# ---------------------------------------------------------------------
$ws.Range("D31").NumberFormat = "#,##0"
$ws.Range("D31").Value = 1
$ws.Range("E31").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E31").Value = -100
$ws.Range("G31").NumberFormat = "#,##0"
$ws.Range("G31").Value = 1
$ws.Range("H31").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H31").Value = -100
$ws.Range("J31").Value = 5
$ws.Range("K31").Value = 0

# ---------------------------------------------------------------------
# Row 33 - Traffic Fatalities (week-to-date count dropped back to zero)
# ---------------------------------------------------------------------
$ws.Range("C14").Copy()
$ws.Range("C33").PasteSpecial(-4122)          # xlPasteFormats - adopt the "text zero" style
$ws.Range("C33").Formula = "=""0"""
$ws.Range("C33").Copy()
$ws.Range("C33").PasteSpecial(-4163)          # xlPasteValues - freeze the literal text "0"
$excel.CutCopyMode = $false
